# Excel COM constants used below
$xlContinuous     = 1
$xlThin           = 2
$xlCenter         = -4108
$xlTop            = -4160
$xlPasteFormats   = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Formatting for the "header" cells B1 and A2 -----------------------
# Bold font, thin box border all around, centered horizontally, top-aligned
# vertically.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = $xlCenter
$r1.VerticalAlignment = $xlTop
$r1.Borders.LineStyle = $xlContinuous
$r1.Borders.Weight = $xlThin

# Clone the exact same formatting onto A2 by copying B1's format over --
# applying the same sequence of property writes to A2 directly can produce
# a second, slightly different style entry in this engine, whereas a
# copy/paste-special of formats reuses B1's already-built style.
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
